# Rename the existing sheet "Feuille1" -> "hoho"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "hoho"

# Add a new sheet "haha" right after "hoho"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "haha"

# Fill in the data for the new "haha" sheet
$ws2.Range("A1").Value = "toto"
$ws2.Range("B1").Value = "titi"
$ws2.Range("C1").Value = "tata"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 3

$ws2.Range("A3").Value = "poisson"
$ws2.Range("B3").Value = "mercredi"
$ws2.Range("C3").Value = "haha"

$ws2.Range("A4").Value = 75
$ws2.Range("B4").Value = 89
$ws2.Range("C4").Value = 108

# Update selections to match the target state
$ws1.Range("G21").Select()
$ws2.Range("C8").Select()

# Make "hoho" the active sheet/tab, as in the target (tabSelected=true on sheet1)
$ws1.Activate()
